$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.533.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.061.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.54%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "386.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("E7").Value = "  -0.59%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -1.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0861"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.548.67"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.055.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.06%  "
$ws.Range("E17").Value = "  -2.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "51.602.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("E20").Value = "  +2.37%  "
$ws.Range("E21").Value = "  -1.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0968"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("E25").Value = "  -2.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.23%  "
$ws.Range("E28").Value = "  -1.96%  "
$ws.Range("E29").Value = "  +1.74%  "
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("E31").Value = "  -1.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.03%  "
$ws.Range("E38").Value = "  +2.07%  "
$ws.Range("E39").Value = "  +8.36%  "
$ws.Range("E40").Value = "  +0.66%  "
$ws.Range("E41").Value = "  +1.57%  "
$ws.Range("E42").Value = "  +0.38%  "
$ws.Range("E43").Value = "  -0.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "125.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("E45").Value = "  +1.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.48%  "
$ws.Range("E48").Value = "  +1.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.034.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.362.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.21%  "
$ws.Range("E51").Value = "  +7.00%  "
